# TrialsSetup.xlsx update (2026-02-16 12:00)
# Source data refresh updated the "Days remaining" figures for two trials:
#   - REJOICE (MK-5909-003) : 6  -> 3   (row 6, column B)
#   - ALPINE                : 26 -> 23  (row 8, column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = 3
$ws.Range("B8").Value = 23
